$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 68.392882
$ws.Range("H2").Value = 205.178646
$ws.Range("I2").Value = 0.3817002623156464
$ws.Range("J2").Value = 0.3817002623156463
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 6.437867666666667
$ws.Range("N2").Value = 19.313603
$ws.Range("O2").Value = 0.4097687230856996
$ws.Range("P2").Value = 0.4097687230856996
$ws.Range("Q2").Value = 440.3043236579487
$ws.Range("R2").Value = 3962.738912921538
$ws.Range("S2").Value = 0.156408829090559
$ws.Range("T2").Value = 0.1564088290905589

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 68.392882
$ws.Range("H3").Value = 205.178646
$ws.Range("I3").Value = 0.3817002623156464
$ws.Range("J3").Value = 0.3817002623156463
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 7.162274333333333
$ws.Range("N3").Value = 21.486823
$ws.Range("O3").Value = 0.4558770325701756
$ws.Range("P3").Value = 0.4558770325701756
$ws.Range("Q3").Value = 489.8485833312953
$ws.Range("R3").Value = 4408.637249981658
$ws.Range("S3").Value = 0.1740083829157145
$ws.Range("T3").Value = 0.1740083829157145

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 68.392882
$ws.Range("H4").Value = 205.178646
$ws.Range("I4").Value = 0.3817002623156464
$ws.Range("J4").Value = 0.3817002623156463
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 2.110836666666667
$ws.Range("N4").Value = 6.33251
$ws.Range("O4").Value = 0.1343542443441249
$ws.Range("P4").Value = 0.1343542443441249
$ws.Range("Q4").Value = 144.3662030646067
$ws.Range("R4").Value = 1299.29582758146
$ws.Range("S4").Value = 0.05128305030937291
$ws.Range("T4").Value = 0.05128305030937291

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 60.20577233333334
$ws.Range("H5").Value = 180.617317
$ws.Range("I5").Value = 0.3360080526004068
$ws.Range("J5").Value = 0.3360080526004068
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 6.437867666666667
$ws.Range("N5").Value = 19.313603
$ws.Range("O5").Value = 0.4097687230856996
$ws.Range("P5").Value = 0.4097687230856996
$ws.Range("Q5").Value = 387.5967950514612
$ws.Range("R5").Value = 3488.371155463151
$ws.Range("S5").Value = 0.1376855906605813
$ws.Range("T5").Value = 0.1376855906605813

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 60.20577233333334
$ws.Range("H6").Value = 180.617317
$ws.Range("I6").Value = 0.3360080526004068
$ws.Range("J6").Value = 0.3360080526004068
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 7.162274333333333
$ws.Range("N6").Value = 21.486823
$ws.Range("O6").Value = 0.4558770325701756
$ws.Range("P6").Value = 0.4558770325701756
$ws.Range("Q6").Value = 431.2102579015435
$ws.Range("R6").Value = 3880.892321113891
$ws.Range("S6").Value = 0.1531783539391569
$ws.Range("T6").Value = 0.1531783539391569

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 60.20577233333334
$ws.Range("H7").Value = 180.617317
$ws.Range("I7").Value = 0.3360080526004068
$ws.Range("J7").Value = 0.3360080526004068
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 2.110836666666667
$ws.Range("N7").Value = 6.33251
$ws.Range("O7").Value = 0.1343542443441249
$ws.Range("P7").Value = 0.1343542443441249
$ws.Range("Q7").Value = 127.0845517861856
$ws.Range("R7").Value = 1143.76096607567
$ws.Range("S7").Value = 0.04514410800066863
$ws.Range("T7").Value = 0.04514410800066863

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 50.58089766666667
$ws.Range("H8").Value = 151.742693
$ws.Range("I8").Value = 0.2822916850839468
$ws.Range("J8").Value = 0.2822916850839468
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 6.437867666666667
$ws.Range("N8").Value = 19.313603
$ws.Range("O8").Value = 0.4097687230856996
$ws.Range("P8").Value = 0.4097687230856996
$ws.Range("Q8").Value = 325.6331256392087
$ws.Range("R8").Value = 2930.698130752879
$ws.Range("S8").Value = 0.1156743033345593
$ws.Range("T8").Value = 0.1156743033345593

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 50.58089766666667
$ws.Range("H9").Value = 151.742693
$ws.Range("I9").Value = 0.2822916850839468
$ws.Range("J9").Value = 0.2822916850839468
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 7.162274333333333
$ws.Range("N9").Value = 21.486823
$ws.Range("O9").Value = 0.4558770325701756
$ws.Range("P9").Value = 0.4558770325701756
$ws.Range("Q9").Value = 362.2742651149265
$ws.Range("R9").Value = 3260.468386034339
$ws.Range("S9").Value = 0.1286902957153041
$ws.Range("T9").Value = 0.1286902957153041

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 50.58089766666667
$ws.Range("H10").Value = 151.742693
$ws.Range("I10").Value = 0.2822916850839468
$ws.Range("J10").Value = 0.2822916850839468
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 2.110836666666667
$ws.Range("N10").Value = 6.33251
$ws.Range("O10").Value = 0.1343542443441249
$ws.Range("P10").Value = 0.1343542443441249
$ws.Range("Q10").Value = 106.7680134277144
$ws.Range("R10").Value = 960.91212084943
$ws.Range("S10").Value = 0.03792708603408333
$ws.Range("T10").Value = 0.03792708603408333
